$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Fecha(D), Calidad(L), Volumen(M), Precio minimo(N), Precio maximo(O), Precio promedio ponderado(P), Origen(R), Precio $/Kg(S)
# This workbook edit re-orders (permutes) the weekly price records across rows 2-32,
# so each row's date/quality/volume/price/origin values are rewritten to their new row position.

$ws.Cells.Item(2, 4).Value = 44343
$ws.Cells.Item(2, 12).Value = "Especial"
$ws.Cells.Item(2, 13).Value = 47
$ws.Cells.Item(2, 14).Value = 10000
$ws.Cells.Item(2, 15).Value = 10000
$ws.Cells.Item(2, 16).Value = 10000
$ws.Cells.Item(2, 18).Value = "Región Metropolitana"
$ws.Cells.Item(2, 19).Value = 1000

$ws.Cells.Item(3, 4).Value = 44343
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 50
$ws.Cells.Item(3, 14).Value = 9000
$ws.Cells.Item(3, 15).Value = 9000
$ws.Cells.Item(3, 16).Value = 9000
$ws.Cells.Item(3, 18).Value = "Región Metropolitana"
$ws.Cells.Item(3, 19).Value = 900

$ws.Cells.Item(4, 4).Value = 44343
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 58
$ws.Cells.Item(4, 14).Value = 8000
$ws.Cells.Item(4, 15).Value = 8000
$ws.Cells.Item(4, 16).Value = 8000
$ws.Cells.Item(4, 18).Value = "Región Metropolitana"
$ws.Cells.Item(4, 19).Value = 800

$ws.Cells.Item(5, 4).Value = 44323
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 60
$ws.Cells.Item(5, 14).Value = 10000
$ws.Cells.Item(5, 15).Value = 10000
$ws.Cells.Item(5, 16).Value = 10000
$ws.Cells.Item(5, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(5, 19).Value = 1000

$ws.Cells.Item(6, 4).Value = 44323
$ws.Cells.Item(6, 12).Value = "Segunda"
$ws.Cells.Item(6, 13).Value = 50
$ws.Cells.Item(6, 14).Value = 9000
$ws.Cells.Item(6, 15).Value = 9000
$ws.Cells.Item(6, 16).Value = 9000
$ws.Cells.Item(6, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(6, 19).Value = 900

$ws.Cells.Item(7, 4).Value = 44322
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 56
$ws.Cells.Item(7, 14).Value = 10000
$ws.Cells.Item(7, 15).Value = 10000
$ws.Cells.Item(7, 16).Value = 10000
$ws.Cells.Item(7, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(7, 19).Value = 1000

$ws.Cells.Item(8, 4).Value = 44322
$ws.Cells.Item(8, 12).Value = "Segunda"
$ws.Cells.Item(8, 13).Value = 40
$ws.Cells.Item(8, 14).Value = 8000
$ws.Cells.Item(8, 15).Value = 8000
$ws.Cells.Item(8, 16).Value = 8000
$ws.Cells.Item(8, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(8, 19).Value = 800

$ws.Cells.Item(9, 4).Value = 44326
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 65
$ws.Cells.Item(9, 14).Value = 10000
$ws.Cells.Item(9, 15).Value = 10000
$ws.Cells.Item(9, 16).Value = 10000
$ws.Cells.Item(9, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(9, 19).Value = 1000

$ws.Cells.Item(10, 4).Value = 44326
$ws.Cells.Item(10, 12).Value = "Segunda"
$ws.Cells.Item(10, 13).Value = 67
$ws.Cells.Item(10, 14).Value = 8000
$ws.Cells.Item(10, 15).Value = 8000
$ws.Cells.Item(10, 16).Value = 8000
$ws.Cells.Item(10, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(10, 19).Value = 800

$ws.Cells.Item(11, 4).Value = 44319
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 68
$ws.Cells.Item(11, 14).Value = 10000
$ws.Cells.Item(11, 15).Value = 10000
$ws.Cells.Item(11, 16).Value = 10000
$ws.Cells.Item(11, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(11, 19).Value = 1000

$ws.Cells.Item(12, 4).Value = 44319
$ws.Cells.Item(12, 12).Value = "Segunda"
$ws.Cells.Item(12, 13).Value = 57
$ws.Cells.Item(12, 14).Value = 8000
$ws.Cells.Item(12, 15).Value = 8000
$ws.Cells.Item(12, 16).Value = 8000
$ws.Cells.Item(12, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(12, 19).Value = 800

$ws.Cells.Item(13, 4).Value = 44307
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 40
$ws.Cells.Item(13, 14).Value = 10000
$ws.Cells.Item(13, 15).Value = 10000
$ws.Cells.Item(13, 16).Value = 10000
$ws.Cells.Item(13, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(13, 19).Value = 1000

$ws.Cells.Item(14, 4).Value = 44312
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 48
$ws.Cells.Item(14, 14).Value = 10000
$ws.Cells.Item(14, 15).Value = 10000
$ws.Cells.Item(14, 16).Value = 10000
$ws.Cells.Item(14, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(14, 19).Value = 1000

$ws.Cells.Item(15, 4).Value = 44301
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 45
$ws.Cells.Item(15, 14).Value = 10000
$ws.Cells.Item(15, 15).Value = 10000
$ws.Cells.Item(15, 16).Value = 10000
$ws.Cells.Item(15, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(15, 19).Value = 1000

$ws.Cells.Item(16, 4).Value = 44333
$ws.Cells.Item(16, 12).Value = "Especial"
$ws.Cells.Item(16, 13).Value = 58
$ws.Cells.Item(16, 14).Value = 10000
$ws.Cells.Item(16, 15).Value = 10000
$ws.Cells.Item(16, 16).Value = 10000
$ws.Cells.Item(16, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(16, 19).Value = 1000

$ws.Cells.Item(17, 4).Value = 44333
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 65
$ws.Cells.Item(17, 14).Value = 9000
$ws.Cells.Item(17, 15).Value = 9000
$ws.Cells.Item(17, 16).Value = 9000
$ws.Cells.Item(17, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(17, 19).Value = 900

$ws.Cells.Item(18, 4).Value = 44333
$ws.Cells.Item(18, 12).Value = "Segunda"
$ws.Cells.Item(18, 13).Value = 60
$ws.Cells.Item(18, 14).Value = 8000
$ws.Cells.Item(18, 15).Value = 8000
$ws.Cells.Item(18, 16).Value = 8000
$ws.Cells.Item(18, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(18, 19).Value = 800

$ws.Cells.Item(19, 4).Value = 44309
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 45
$ws.Cells.Item(19, 14).Value = 10000
$ws.Cells.Item(19, 15).Value = 10000
$ws.Cells.Item(19, 16).Value = 10000
$ws.Cells.Item(19, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(19, 19).Value = 1000

$ws.Cells.Item(20, 4).Value = 44308
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 45
$ws.Cells.Item(20, 14).Value = 10000
$ws.Cells.Item(20, 15).Value = 10000
$ws.Cells.Item(20, 16).Value = 10000
$ws.Cells.Item(20, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(20, 19).Value = 1000

$ws.Cells.Item(21, 4).Value = 44308
$ws.Cells.Item(21, 12).Value = "Segunda"
$ws.Cells.Item(21, 13).Value = 48
$ws.Cells.Item(21, 14).Value = 8000
$ws.Cells.Item(21, 15).Value = 8000
$ws.Cells.Item(21, 16).Value = 8000
$ws.Cells.Item(21, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(21, 19).Value = 800

$ws.Cells.Item(22, 4).Value = 44302
$ws.Cells.Item(22, 12).Value = "Primera"
$ws.Cells.Item(22, 13).Value = 45
$ws.Cells.Item(22, 14).Value = 10000
$ws.Cells.Item(22, 15).Value = 10000
$ws.Cells.Item(22, 16).Value = 10000
$ws.Cells.Item(22, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(22, 19).Value = 1000

$ws.Cells.Item(23, 4).Value = 44306
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 45
$ws.Cells.Item(23, 14).Value = 10000
$ws.Cells.Item(23, 15).Value = 10000
$ws.Cells.Item(23, 16).Value = 10000
$ws.Cells.Item(23, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(23, 19).Value = 1000

$ws.Cells.Item(24, 4).Value = 44699
$ws.Cells.Item(24, 12).Value = "Especial"
$ws.Cells.Item(24, 13).Value = 56
$ws.Cells.Item(24, 14).Value = 12000
$ws.Cells.Item(24, 15).Value = 12000
$ws.Cells.Item(24, 16).Value = 12000
$ws.Cells.Item(24, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(24, 19).Value = 1200

$ws.Cells.Item(25, 4).Value = 44699
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 60
$ws.Cells.Item(25, 14).Value = 10000
$ws.Cells.Item(25, 15).Value = 10000
$ws.Cells.Item(25, 16).Value = 10000
$ws.Cells.Item(25, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(25, 19).Value = 1000

$ws.Cells.Item(26, 4).Value = 44321
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 58
$ws.Cells.Item(26, 14).Value = 9000
$ws.Cells.Item(26, 15).Value = 9000
$ws.Cells.Item(26, 16).Value = 9000
$ws.Cells.Item(26, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(26, 19).Value = 900

$ws.Cells.Item(27, 4).Value = 44329
$ws.Cells.Item(27, 12).Value = "Primera"
$ws.Cells.Item(27, 13).Value = 56
$ws.Cells.Item(27, 14).Value = 9000
$ws.Cells.Item(27, 15).Value = 9000
$ws.Cells.Item(27, 16).Value = 9000
$ws.Cells.Item(27, 18).Value = "Región Metropolitana"
$ws.Cells.Item(27, 19).Value = 900

$ws.Cells.Item(28, 4).Value = 44329
$ws.Cells.Item(28, 12).Value = "Segunda"
$ws.Cells.Item(28, 13).Value = 50
$ws.Cells.Item(28, 14).Value = 8000
$ws.Cells.Item(28, 15).Value = 8000
$ws.Cells.Item(28, 16).Value = 8000
$ws.Cells.Item(28, 18).Value = "Región Metropolitana"
$ws.Cells.Item(28, 19).Value = 800

$ws.Cells.Item(29, 4).Value = 44315
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 45
$ws.Cells.Item(29, 14).Value = 10000
$ws.Cells.Item(29, 15).Value = 10000
$ws.Cells.Item(29, 16).Value = 10000
$ws.Cells.Item(29, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(29, 19).Value = 1000

$ws.Cells.Item(30, 4).Value = 44328
$ws.Cells.Item(30, 12).Value = "Primera"
$ws.Cells.Item(30, 13).Value = 45
$ws.Cells.Item(30, 14).Value = 8000
$ws.Cells.Item(30, 15).Value = 8000
$ws.Cells.Item(30, 16).Value = 8000
$ws.Cells.Item(30, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(30, 19).Value = 800

$ws.Cells.Item(31, 4).Value = 44328
$ws.Cells.Item(31, 12).Value = "Segunda"
$ws.Cells.Item(31, 13).Value = 48
$ws.Cells.Item(31, 14).Value = 7000
$ws.Cells.Item(31, 15).Value = 7000
$ws.Cells.Item(31, 16).Value = 7000
$ws.Cells.Item(31, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(31, 19).Value = 700

$ws.Cells.Item(32, 4).Value = 44314
$ws.Cells.Item(32, 12).Value = "Primera"
$ws.Cells.Item(32, 13).Value = 47
$ws.Cells.Item(32, 14).Value = 9000
$ws.Cells.Item(32, 15).Value = 9000
$ws.Cells.Item(32, 16).Value = 9000
$ws.Cells.Item(32, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(32, 19).Value = 900
